$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the additional precondition line to B9, and grow the row so the
# now-two-line wrapped text is fully visible.
$ws.Range("B9").Value = "Que el actor tenga los permisos necesarios para ver el registro." + [char]10 + "Que existan campañas vigentes."
$ws.Rows.Item(9).RowHeight = 25.5

# Update version number (B5): 0001 -> 0002
$ws.Range("B5").Value = "0002"

# Move the active selection to B6 (matches the saved view state).
$ws.Range("B6").Select()
